{"js": "// Find the \"Electives\" heading paragraph and the list-of-electives paragraph\n// right after it, then:\n//  1. Turn \"Electives\" into a proper \"Electives:\" Heading 4 (matching the\n//     other section headings like \"Software Development Projects:\"), with a\n//     bookmark named \"electives\" around it.\n//  2. Promote the paragraph describing the electives (the \"Machine Learning\n//     & Data Mining, ...\" line) from Body Text to First Paragraph, matching\n//     the style used for the descriptive line right under other headings.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nlet electivesHeading = null;\nlet electivesBody = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === \"Electives\") {\n    electivesHeading = para;\n    if (i + 1 < paragraphs.items.length) {\n      electivesBody = paragraphs.items[i + 1];\n    }\n    break;\n  }\n}\n\nif (!electivesHeading) {\n  throw new Error('Could not find the \"Electives\" paragraph.');\n}\n\n// Replace the run text (dropping the bold formatting) and restyle as Heading 4.\nconst headingRange = electivesHeading.getRange();\nheadingRange.insertText(\"Electives:\", \"Replace\");\nelectivesHeading.style = \"Heading 4\";\n\n// Mark it with a bookmark the same way the other headings are bookmarked\n// (a zero-width bookmark sitting right at the start of the paragraph, same\n// as \"software-development-projects\", \"skills\", etc.).\nelectivesHeading.getRange(\"Start\").insertBookmark(\"electives\");\n\n// The elective-course list paragraph becomes a \"First Paragraph\" (same style\n// used for the descriptive paragraph right below other Heading 4 titles).\nif (electivesBody) {\n  electivesBody.style = \"First Paragraph\";\n}\n\nawait context.sync();\n", "ps1": "# Find the \"Electives\" heading paragraph and the list-of-electives paragraph\n# right after it, then:\n#  1. Turn \"Electives\" into a proper \"Electives:\" Heading 4 (matching the\n#     other section headings like \"Software Development Projects:\"), with a\n#     bookmark named \"electives\" around it.\n#  2. Promote the paragraph describing the electives (the \"Machine Learning\n#     & Data Mining, ...\" line) from Body Text to First Paragraph, matching\n#     the style used for the descriptive line right under other headings.\n\n$d = $word.ActiveDocument\n\n$electivesIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text.Trim() -eq \"Electives\") {\n        $electivesIndex = $i\n        break\n    }\n}\n\nif ($electivesIndex -eq -1) {\n    throw \"Could not find the 'Electives' paragraph.\"\n}\n\n$heading = $d.Paragraphs.Item($electivesIndex)\n\n# Replace the run text (dropping the bold formatting) and restyle as Heading 4.\n$heading.Range.Text = \"Electives:\"\n$heading.Style = \"Heading 4\"\n\n# Mark it with a bookmark the same way the other headings are bookmarked\n# (a zero-width bookmark sitting right at the start of the paragraph, same\n# as \"software-development-projects\", \"skills\", etc.).\n$bmRange = $heading.Range.Duplicate\n$bmRange.Collapse(1)\n$d.Bookmarks.Add(\"electives\", $bmRange)\n\n# The elective-course list paragraph becomes a \"First Paragraph\" (same style\n# used for the descriptive paragraph right below other Heading 4 titles).\n$bodyIndex = $electivesIndex + 1\nif ($bodyIndex -le $d.Paragraphs.Count) {\n    $electivesBody = $d.Paragraphs.Item($bodyIndex)\n    $electivesBody.Style = \"First Paragraph\"\n}\n"}
